$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4169947665850486
$ws.Range("D2").Value = 0.1383540777476355
$ws.Range("E2").Value = 0.1767078081192253
$ws.Range("F2").Value = 2.21055258455722
$ws.Range("G2").Value = 0.002520715835772573
$ws.Range("J2").Value = 0.2511392977561329
$ws.Range("K2").Value = 1.965132290912948
$ws.Range("L2").Value = 0.1562234536183666
$ws.Range("N2").Value = 1.492239119467051
$ws.Range("O2").Value = 5.840205772284008
$ws.Range("C3").Value = 0.412471168722476
$ws.Range("D3").Value = 0.1358399787789466
$ws.Range("E3").Value = 0.1756992814456986
$ws.Range("F3").Value = 2.215540238303277
$ws.Range("G3").Value = 0.00252411860680124
$ws.Range("J3").Value = 0.2511006519643786
$ws.Range("K3").Value = 1.834975815562188
$ws.Range("L3").Value = 0.1559134350462941
$ws.Range("N3").Value = 1.500826239000787
$ws.Range("O3").Value = 5.864365450898248
$ws.Range("C4").Value = 0.4098828058754123
$ws.Range("D4").Value = 0.1343390114345553
$ws.Range("E4").Value = 0.1751481770980625
$ws.Range("F4").Value = 2.219756885850146
$ws.Range("G4").Value = 0.00252632007684012
$ws.Range("J4").Value = 0.2511831465795922
$ws.Range("K4").Value = 1.75547353548518
$ws.Range("L4").Value = 0.1557718132252575
$ws.Range("N4").Value = 1.5066498965192
$ws.Range("O4").Value = 5.882379208610587
$ws.Range("C5").Value = 0.4088756748296447
$ws.Range("D5").Value = 0.133738159557268
$ws.Range("E5").Value = 0.1749407698995675
$ws.Range("F5").Value = 2.221765386482893
$ws.Range("G5").Value = 0.002527245482265784
$ws.Range("J5").Value = 0.2512434995867281
$ws.Range("K5").Value = 1.723181888151004
$ws.Range("L5").Value = 0.1557263886807014
$ws.Range("N5").Value = 1.509161892703247
$ws.Range("O5").Value = 5.890519048516893
$ws.Range("C6").Value = 0.4087113218426026
$ws.Range("D6").Value = 0.1336390431130638
$ws.Range("E6").Value = 0.1749073684141997
$ws.Range("F6").Value = 2.222116421330526
$ws.Range("G6").Value = 0.002527400856304225
$ws.Range("J6").Value = 0.251255136677905
$ws.Range("K6").Value = 1.7178263564264
$ws.Range("L6").Value = 0.1557195890993093
$ws.Range("N6").Value = 1.509587398100138
$ws.Range("O6").Value = 5.89191891254589
$ws.Range("C7").Value = 0.4098690303141126
$ws.Range("D7").Value = 0.1343308642984908
$ws.Range("E7").Value = 0.1751453103496452
$ws.Range("F7").Value = 2.219782798327714
$ws.Range("G7").Value = 0.00252633244236487
$ws.Range("J7").Value = 0.2511838522354921
$ws.Range("K7").Value = 1.755037606483938
$ws.Range("L7").Value = 0.1557711508173973
$ws.Range("N7").Value = 1.50668321185195
$ws.Range("O7").Value = 5.882485750478708
$ws.Range("C8").Value = 0.4153958221059639
$ws.Range("D8").Value = 0.1374783934446313
$ws.Range("E8").Value = 0.1763459494482262
$ws.Range("F8").Value = 2.212032735030448
$ws.Range("G8").Value = 0.002521865885565756
$ws.Range("J8").Value = 0.2511039404246134
$ws.Range("K8").Value = 1.920169496019923
$ws.Range("L8").Value = 0.1561064649591977
$ws.Range("N8").Value = 1.495085701533903
$ws.Range("O8").Value = 5.847875816759711
$ws.Range("C9").Value = 0.4277318799456395
$ws.Range("D9").Value = 0.1439869449418012
$ws.Range("E9").Value = 0.1792395929264465
$ws.Range("F9").Value = 2.205997684546745
$ws.Range("G9").Value = 0.002513992932883863
$ws.Range("J9").Value = 0.2517892415458007
$ws.Range("K9").Value = 2.247209764318427
$ws.Range("L9").Value = 0.1571493478542081
$ws.Range("N9").Value = 1.476706857979295
$ws.Range("O9").Value = 5.805259942585735
$ws.Range("C10").Value = 0.4377062686563988
$ws.Range("D10").Value = 0.1489709560464263
$ws.Range("E10").Value = 0.1816927744458638
$ws.Range("F10").Value = 2.207159442769125
$ws.Range("G10").Value = 0.002508743245777633
$ws.Range("J10").Value = 0.2528053407358115
$ws.Range("K10").Value = 2.489378389368653
$ws.Range("L10").Value = 0.1581489155040359
$ws.Range("N10").Value = 1.465852841362192
$ws.Range("O10").Value = 5.789387586186024
$ws.Range("C11").Value = 0.4424413604361632
$ws.Range("D11").Value = 0.151281657972433
$ws.Range("E11").Value = 0.1828795599944861
$ws.Range("F11").Value = 2.208905400189181
$ws.Range("G11").Value = 0.002506469927628263
$ws.Range("J11").Value = 0.2533787833892944
$ws.Range("K11").Value = 2.599945067130818
$ws.Range("L11").Value = 0.1586540067806297
$ws.Range("N11").Value = 1.461488033156087
$ws.Range("O11").Value = 5.785528599439459
$ws.Range("C12").Value = 0.4442627841852698
$ws.Range("D12").Value = 0.1521628497895193
$ws.Range("E12").Value = 0.1833391168356826
$ws.Range("F12").Value = 2.209741763204264
$ws.Range("G12").Value = 0.002505625498711206
$ws.Range("J12").Value = 0.2536119072922816
$ws.Range("K12").Value = 2.641870096329001
$ws.Range("L12").Value = 0.158852486809451
$ws.Range("N12").Value = 1.459917379960359
$ws.Range("O12").Value = 5.784551365314996
$ws.Range("C13").Value = 0.4438692485181264
$ws.Range("D13").Value = 0.1519727956412282
$ws.Range("E13").Value = 0.1832396923486499
$ws.Range("F13").Value = 2.20955384289627
$ws.Range("G13").Value = 0.002505806632220418
$ws.Range("J13").Value = 0.2535609897545612
$ws.Range("K13").Value = 2.632838341779348
$ws.Range("L13").Value = 0.1588094202959169
$ws.Range("N13").Value = 1.460251994516014
$ws.Range("O13").Value = 5.784740289166962
$ws.Range("C14").Value = 0.4425906424026493
$ws.Range("D14").Value = 0.1513540307080632
$ws.Range("E14").Value = 0.1829171648646692
$ws.Range("F14").Value = 2.208970696082815
$ws.Range("G14").Value = 0.002506400127118162
$ws.Range("J14").Value = 0.2533976426048028
$ws.Range("K14").Value = 2.603393157726032
$ws.Range("L14").Value = 0.1586701914774054
$ws.Range("N14").Value = 1.461357167919559
$ws.Range("O14").Value = 5.785438496012347
$ws.Range("C15").Value = 0.4418111485022109
$ws.Range("D15").Value = 0.1509758219858242
$ws.Range("E15").Value = 0.182720927799064
$ws.Range("F15").Value = 2.208636323009259
$ws.Range("G15").Value = 0.00250676579676791
$ws.Range("J15").Value = 0.2532996674424126
$ws.Range("K15").Value = 2.585364340091473
$ws.Range("L15").Value = 0.1585858481923523
$ws.Range("N15").Value = 1.462044819602511
$ws.Range("O15").Value = 5.785929230635418
$ws.Range("C16").Value = 0.4374007904415294
$ws.Range("D16").Value = 0.1488208146714527
$ws.Range("E16").Value = 0.181616637145158
$ws.Range("F16").Value = 2.207069850523879
$ws.Range("G16").Value = 0.002508894115732881
$ws.Range("J16").Value = 0.2527701012034953
$ws.Range("K16").Value = 2.48216052426784
$ws.Range("L16").Value = 0.1581169176880834
$ws.Range("N16").Value = 1.466149603457211
$ws.Range("O16").Value = 5.789707481217391
$ws.Range("C17").Value = 0.4347457548241209
$ws.Range("D17").Value = 0.1475098669894521
$ws.Range("E17").Value = 0.1809573031071849
$ws.Range("F17").Value = 2.206420789504747
$ws.Range("G17").Value = 0.002510229114264046
$ws.Range("J17").Value = 0.2524737005561377
$ws.Range("K17").Value = 2.418950023245486
$ws.Range("L17").Value = 0.1578421254049047
$ws.Range("N17").Value = 1.468814335398861
$ws.Range("O17").Value = 5.79288676195506
$ws.Range("C18").Value = 0.4332372610913069
$ws.Range("D18").Value = 0.1467599384069302
$ws.Range("E18").Value = 0.1805847398148401
$ws.Range("F18").Value = 2.206162052290566
$ws.Range("G18").Value = 0.002511007779590228
$ws.Range("J18").Value = 0.2523136878415571
$ws.Range("K18").Value = 2.382631068832723
$ws.Range("L18").Value = 0.1576888163381511
$ws.Range("N18").Value = 1.470400936305836
$ws.Range("O18").Value = 5.795031750671455
$ws.Range("C19").Value = 0.4327297102187799
$ws.Range("D19").Value = 0.1465067307021855
$ws.Range("E19").Value = 0.1804597426892371
$ws.Range("F19").Value = 2.206094124374573
$ws.Range("G19").Value = 0.002511273281211349
$ws.Range("J19").Value = 0.2522613089422734
$ws.Range("K19").Value = 2.370340695385494
$ws.Range("L19").Value = 0.1576377243214182
$ws.Range("N19").Value = 1.470947397330377
$ws.Range("O19").Value = 5.795812318960174
$ws.Range("C20").Value = 0.4350264621578788
$ws.Range("D20").Value = 0.1476489964107657
$ws.Range("E20").Value = 0.1810268004761078
$ws.Range("F20").Value = 2.206478023000003
$ws.Range("G20").Value = 0.002510085883299432
$ws.Range("J20").Value = 0.252504169536266
$ws.Range("K20").Value = 2.425674965913345
$ws.Range("L20").Value = 0.1578708867398291
$ws.Range("N20").Value = 1.468525091004722
$ws.Range("O20").Value = 5.792515575850075
$ws.Range("C21").Value = 0.442965431317333
$ws.Range("D21").Value = 0.1515356098151983
$ws.Range("E21").Value = 0.1830116239149042
$ws.Range("F21").Value = 2.209137224276091
$ws.Range("G21").Value = 0.00250622535806603
$ws.Range("J21").Value = 0.2534451882516251
$ws.Range("K21").Value = 2.612040423705196
$ws.Range("L21").Value = 0.1587108908541239
$ws.Range("N21").Value = 1.461030321873764
$ws.Range("O21").Value = 5.78522027223454
$ws.Range("C22").Value = 0.448319199029271
$ws.Range("D22").Value = 0.154111721362824
$ws.Range("E22").Value = 0.1843679476926177
$ws.Range("F22").Value = 2.211896470423582
$ws.Range("G22").Value = 0.00250379799734984
$ws.Range("J22").Value = 0.2541532963213555
$ws.Range("K22").Value = 2.734165457110407
$ws.Range("L22").Value = 0.1593019097953601
$ws.Range("N22").Value = 1.456611144037581
$ws.Range("O22").Value = 5.783274295173669
$ws.Range("C23").Value = 0.4454467038448797
$ws.Range("D23").Value = 0.1527335322214185
$ws.Range("E23").Value = 0.1836386542348727
$ws.Range("F23").Value = 2.210330313760437
$ws.Range("G23").Value = 0.002505084793639304
$ws.Range("J23").Value = 0.2537668533890809
$ws.Range("K23").Value = 2.66895603928873
$ws.Range("L23").Value = 0.1589826368759404
$ws.Range("N23").Value = 1.458925954929455
$ws.Range("O23").Value = 5.784054456793058
$ws.Range("C24").Value = 0.4348994985440129
$ws.Range("D24").Value = 0.1475860843056864
$ws.Range("E24").Value = 0.1809953604770236
$ws.Range("F24").Value = 2.206451791338452
$ws.Range("G24").Value = 0.002510150603245484
$ws.Range("J24").Value = 0.2524903621394969
$ws.Range("K24").Value = 2.422634551871852
$ws.Range("L24").Value = 0.157857869182962
$ws.Range("N24").Value = 1.46865568826783
$ws.Range("O24").Value = 5.792682401124438
$ws.Range("C25").Value = 0.4242344992283051
$ws.Range("D25").Value = 0.1421905026321042
$ws.Range("E25").Value = 0.1783992157009671
$ws.Range("F25").Value = 2.206648378571131
$ws.Range("G25").Value = 0.002516028501598991
$ws.Range("J25").Value = 0.2515137804882031
$ws.Range("K25").Value = 2.158399475755516
$ws.Range("L25").Value = 0.1568261285801995
$ws.Range("N25").Value = 1.48121285028828
$ws.Range("O25").Value = 5.814080729789566
